$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the magic-def reference strings in column D/E (rows 4-7).
# The underlying numeric ids are all decremented by one (bug fix for
# "null magic def" being off-by-one), e.g. 51018002 -> 51018001 and
# 51018003 -> 51018002.
$ws.Range("D4").Value = "0;4;51018001;2;1;51018002;2;7;51018002"
$ws.Range("E4").Value = "19;4;51018001;17;1;51018002;17;7;51018002"

$ws.Range("D5").Value = "1;2;51018001"
$ws.Range("E5").Value = "9;2;51018001"

$ws.Range("D6").Value = "0;2;51018001"
$ws.Range("E6").Value = "10;2;51018001"

$ws.Range("D7").Value = "8;4;51018001"
$ws.Range("E7").Value = "14;4;51018001"

# Match the author's new cursor/selection position (was F5, now E5).
$ws.Range("E5").Select()
